# Scheduled-runner market-data refresh: updates the static
# currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) on a handful of
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets with freshly pulled
# values. These columns hold plain numbers (no formulas in this workbook),
# so each cell is written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 278.2381
$ws.Range("I28").Value = 171.0625
$ws.Range("J28").Value = 621.2
$ws.Range("K28").Value = 171.0625
$ws.Range("L28").Value = 621.2
$ws.Range("M28").Value = 313.9375
$ws.Range("N28").Value = -1591.2

$ws.Range("H80").Value = 342.41934
$ws.Range("I80").Value = 338.8889
$ws.Range("J80").Value = 347.30768
$ws.Range("K80").Value = 1016.6667
$ws.Range("L80").Value = 1041.92304
$ws.Range("M80").Value = -18.66669999999999
$ws.Range("N80").Value = -3037.92304

$ws.Range("H83").Value = 342.41934
$ws.Range("I83").Value = 338.8889
$ws.Range("J83").Value = 347.30768
$ws.Range("K83").Value = 3050.0001
$ws.Range("L83").Value = 3125.76912
$ws.Range("M83").Value = 1941.9999
$ws.Range("N83").Value = -13109.76912

$ws.Range("H92").Value = 514.5333000000001
$ws.Range("I92").Value = 299.57144
$ws.Range("J92").Value = 1016.1111
$ws.Range("K92").Value = 299.57144
$ws.Range("L92").Value = 1016.1111
$ws.Range("M92").Value = 948.4285600000001
$ws.Range("N92").Value = -3512.1111

$ws.Range("H106").Value = 35715708
$ws.Range("I106").Value = 52632430
$ws.Range("J106").Value = 2620
$ws.Range("K106").Value = 52632430
$ws.Range("L106").Value = 2620
$ws.Range("M106").Value = -52631799
$ws.Range("N106").Value = -3882

$ws.Range("H138").Value = 2008.52
$ws.Range("I138").Value = 933.7692
$ws.Range("J138").Value = 2695.6558
$ws.Range("K138").Value = 2801.3076
$ws.Range("L138").Value = 8086.9674
$ws.Range("M138").Value = 2338.6924
$ws.Range("N138").Value = -18366.9674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5593.71
$ws.Range("I32").Value = 5203.3735
$ws.Range("J32").Value = 7499.4707
$ws.Range("K32").Value = 5203.3735
$ws.Range("L32").Value = 7499.4707
$ws.Range("M32").Value = -4916.3735
$ws.Range("N32").Value = -8073.4707

$ws.Range("H61").Value = 2009.7949
$ws.Range("I61").Value = 1933.0454
$ws.Range("J61").Value = 2109.1177
$ws.Range("K61").Value = 1933.0454
$ws.Range("L61").Value = 2109.1177
$ws.Range("M61").Value = -1721.0454
$ws.Range("N61").Value = -2533.1177

$ws.Range("H74").Value = 15152683
$ws.Range("I74").Value = 16130181
$ws.Range("J74").Value = 1449.5
$ws.Range("K74").Value = 16130181
$ws.Range("L74").Value = 1449.5
$ws.Range("M74").Value = -16129307
$ws.Range("N74").Value = -3197.5

$ws.Range("H77").Value = 15152683
$ws.Range("I77").Value = 16130181
$ws.Range("J77").Value = 1449.5
$ws.Range("K77").Value = 80650905
$ws.Range("L77").Value = 7247.5
$ws.Range("M77").Value = -80646537
$ws.Range("N77").Value = -15983.5

$ws.Range("H136").Value = 2009.7949
$ws.Range("I136").Value = 1933.0454
$ws.Range("J136").Value = 2109.1177
$ws.Range("K136").Value = 5799.1362
$ws.Range("L136").Value = 6327.353099999999
$ws.Range("M136").Value = -3249.1362
$ws.Range("N136").Value = -11427.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8333814
$ws.Range("I107").Value = 10638683
$ws.Range("J107").Value = 828.38464
$ws.Range("K107").Value = 10638683
$ws.Range("L107").Value = 828.38464
$ws.Range("M107").Value = -10636763
$ws.Range("N107").Value = -4668.38464

$ws.Range("H134").Value = 1766631.5
$ws.Range("I134").Value = 1033.6428
$ws.Range("J134").Value = 5297827.5
$ws.Range("K134").Value = 3100.9284
$ws.Range("L134").Value = 15893482.5
$ws.Range("M134").Value = -565.9284000000002
$ws.Range("N134").Value = -15898552.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 669.75
$ws.Range("J10").Value = 800
$ws.Range("L10").Value = 800
$ws.Range("N10").Value = -1078

$ws.Range("H16").Value = 1838.4231
$ws.Range("I16").Value = 857.0714
$ws.Range("J16").Value = 2983.3333
$ws.Range("K16").Value = 857.0714
$ws.Range("L16").Value = 2983.3333
$ws.Range("M16").Value = -570.0714
$ws.Range("N16").Value = -3557.3333

$ws.Range("H36").Value = 2666.6667
$ws.Range("I36").Value = 2666.6667
$ws.Range("K36").Value = 2666.6667
$ws.Range("M36").Value = -2278.6667

$ws.Range("H40").Value = 2666.6667
$ws.Range("I40").Value = 2666.6667
$ws.Range("K40").Value = 2666.6667
$ws.Range("M40").Value = -2506.6667

$ws.Range("H113").Value = 1838.4231
$ws.Range("I113").Value = 857.0714
$ws.Range("J113").Value = 2983.3333
$ws.Range("K113").Value = 857.0714
$ws.Range("L113").Value = 2983.3333
$ws.Range("M113").Value = 1312.9286
$ws.Range("N113").Value = -7323.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 10784855
$ws.Range("I113").Value = 7576271.5
$ws.Range("J113").Value = 16667258
$ws.Range("K113").Value = 22728814.5
$ws.Range("L113").Value = 50001774
$ws.Range("M113").Value = -22726644.5
$ws.Range("N113").Value = -50006114

$ws.Range("H131").Value = 713.99
$ws.Range("I131").Value = 413.6842
$ws.Range("J131").Value = 784.4321
$ws.Range("K131").Value = 1241.0526
$ws.Range("L131").Value = 2353.2963
$ws.Range("M131").Value = 3798.9474
$ws.Range("N131").Value = -12433.2963

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 909
$ws.Range("I102").Value = 887.175
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 887.175
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 734.825
$ws.Range("N102").Value = -4444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2388
$ws.Range("N7").ClearContents()

$ws.Range("I14").Value = 2500
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2328
$ws.Range("N14").ClearContents()

$ws.Range("H122").Value = 7458.5366
$ws.Range("I122").Value = 8503.75
$ws.Range("J122").Value = 3742.2222
$ws.Range("K122").Value = 25511.25
$ws.Range("L122").Value = 11226.6666
$ws.Range("M122").Value = -23061.25
$ws.Range("N122").Value = -16126.6666

$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5030
$ws.Range("N126").ClearContents()
